# --- "Metadata" sheet: update Version / Status / Date / Contact -----------
$wb   = $excel.ActiveWorkbook
$meta = $wb.Worksheets.Item("Metadata")

$meta.Range("B3").Value  = "0.4.0-snapshot-1"              # Version
$meta.Range("B6").Value  = "draft"                         # Status
$meta.Range("B8").Value  = "2024-05-23T12:16:26+00:00"     # Date
$meta.Range("B10").Value = "ANS (https://esante.gouv.fr)"  # Contact

# --- "Elements" sheet: swap the "AK" / "AL" mapping columns ---------------
# (the "RIM Mapping" column and the "Spécification métier vers l'extension
#  ROR AccomodationFamily" column traded places, header + data + width)
$els = $wb.Worksheets.Item("Elements")

for ($r = 1; $r -le 6; $r++) {
    $akCell = $els.Cells.Item($r, 37)
    $alCell = $els.Cells.Item($r, 38)
    $akVal  = $akCell.Value2
    $alVal  = $alCell.Value2

    if ($akVal -ne $alVal) {
        if ($alVal -eq $null) {
            $akCell.ClearContents()
        } else {
            $akCell.Value = $alVal
        }

        if ($akVal -eq $null) {
            $alCell.ClearContents()
        } else {
            $alCell.Value = $akVal
        }
    }
}

# Column widths also swap along with the content (AK becomes the wide
# "Spécification métier" column, AL becomes the narrower "RIM Mapping" one).
$els.Columns.Item(37).ColumnWidth = 75.66666666666667   # -> stored width ~76.53
$els.Columns.Item(38).ColumnWidth = 24.166666666666668  # -> stored width ~24.98
